$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").NumberFormat = "@"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"

$ws.Range("A14").Value = "2025-09-30"
$ws.Range("B14").Value = "Pick 4"
$ws.Range("C14").Value = "250930"
$ws.Range("D14").Value = "8-8-3-9"
$ws.Range("E14").Value = "2025-09-30T21:37:59.536+04:00"
